# Update the ACTION cell (F6) from "autogen" to "autogen == $param"
# and move the active selection from F5 to F7, matching the authored
# decision-table edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F6").Value = "autogen == `$param"

$ws.Activate()
$ws.Range("F7").Select()
